# Corrected excel sheets for application fix issues
#
# This script reproduces, via Excel COM-interop semantics, the edits that
# were made to the workbook:
#   - Summary sheet: a few corrected figures + new selection
#   - Repayment schedule sheet: a new "Over Due" (O) column of zeros, a new
#     repayment schedule row (row 8) for a follow-up installment, and a
#     couple of corrected Interest/Paid figures
#   - Transactions sheet: corrected transaction ids, and this sheet becomes
#     the active tab/selection

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("F2").Value = 0
$summary.Range("A3").Value = 211.3
$summary.Range("E3").Value = 114.31

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$repay = $wb.Worksheets.Item("Repayment schedule")

# New "Over Due" column O (between N "Late" and P "Outstanding"): copy the
# formatting from column N, then fill in the values (row 2 stays blank like
# N2, rows 3-7 get 0).
$repay.Range("N2:N7").Copy()
$repay.Range("O2:O7").PasteSpecial(-4122)
$repay.Range("O3:O7").Value = 0

# Corrected figures on existing rows
$repay.Range("H5").Value = 31.38
$repay.Range("K5").Value = 940.47
$repay.Range("P5").Value = 940.47

$repay.Range("H6").Value = 23.16
$repay.Range("K6").Value = 932.25
$repay.Range("P6").Value = 932.25

$repay.Range("F7").Value = 909.09
$repay.Range("H7").Value = 13.45
$repay.Range("K7").Value = 922.54
$repay.Range("P7").Value = 922.54

# New row 8 - another installment. Copy formatting from row 7 first so the
# new row matches the existing styles (including the new O column), then
# set its values.
$repay.Range("A7:P7").Copy()
$repay.Range("A8:P8").PasteSpecial(-4122)

$repay.Range("A8").Value = 6
$repay.Range("B8").Value = 31
$repay.Range("C8").Value = 42217
$repay.Range("D8").Value = ""
$repay.Range("E8").Value = ""
$repay.Range("F8").Value = 454.55
$repay.Range("G8").Value = 0
$repay.Range("H8").Value = 4.63
$repay.Range("I8").Value = 0
$repay.Range("J8").Value = 0
$repay.Range("K8").Value = 459.18
$repay.Range("L8").Value = 0
$repay.Range("M8").Value = 0
$repay.Range("N8").Value = 0
$repay.Range("O8").Value = 0
$repay.Range("P8").Value = 459.18

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$transactions = $wb.Worksheets.Item("Transactions")

$transactions.Range("A2").Value = 98
$transactions.Range("A3").Value = 96

# ---------------------------------------------------------------------------
# Selections / active tab
# ---------------------------------------------------------------------------
# Leave a plain selection sitting on the Summary sheet...
$summary.Range("D4").Select()

# ...and on the Repayment schedule sheet (now one row further down, since a
# row was appended)...
$repay.Range("A9:XFD9").Select()

# ...and finally land on the Transactions sheet, which becomes the active
# tab.
$transactions.Activate()
$transactions.Range("A2:L3").Select()
